$wb = $excel.ActiveWorkbook

# Sheet 1: AMC Aerospace Solutions Divisi
$ws = $wb.Worksheets.Item(1)
$ws.Range("D2").Value = 0.0634
$ws.Range("D3").Value = 0.0634
$ws.Range("D4").Value = 0.0634
$ws.Range("F4").Value = 0.0075
$ws.Range("G4").Value = 0.0114
$ws.Range("H4").Value = 0.019
$ws.Range("I4").Value = 0.0378
$ws.Range("J4").Value = 0.0114
$ws.Range("L4").Value = 0.0073
$ws.Range("M4").Value = 0.0224
$ws.Range("N4").Value = 0.0036
$ws.Range("O4").Value = 0.00905833333333333
$ws.Range("P4").Value = 0.00905833333333333
$ws.Range("Q4").Value = 0.027175
$ws.Range("R4").Value = 0.00905833333333333
$ws.Range("S4").Value = 0.00905833333333333
$ws.Range("T4").Value = 0.00905833333333333
$ws.Range("U4").Value = 0.027175
$ws.Range("V4").Value = 0.1087
$ws.Range("D5").Value = 0.233333333333333
$ws.Range("D6").Value = 0.233333333333333
$ws.Range("D7").Value = 0.233333333333333
$ws.Range("N7").Value = 0
$ws.Range("O7").Value = 0.233333333333333
$ws.Range("P7").Value = 0.233333333333333
$ws.Range("Q7").Value = 0.233333333333333
$ws.Range("R7").Value = 0.233333333333333
$ws.Range("S7").Value = 0.233333333333333
$ws.Range("T7").Value = 0.233333333333333
$ws.Range("U7").Value = 0.233333333333333
$ws.Range("V7").Value = 0.233333333333333

# Sheet 2: AMC Autonomous Mobile Solution
$ws = $wb.Worksheets.Item(2)
$ws.Range("D2").Value = 0.008
$ws.Range("D3").Value = 0.008
$ws.Range("D4").Value = 0.008
$ws.Range("N4").Value = 0
$ws.Range("O4").Value = 0.00114166666666667
$ws.Range("P4").Value = 0.00114166666666667
$ws.Range("Q4").Value = 0.003425
$ws.Range("R4").Value = 0.00114166666666667
$ws.Range("S4").Value = 0.00114166666666667
$ws.Range("T4").Value = 0.00114166666666667
$ws.Range("U4").Value = 0.003425
$ws.Range("V4").Value = 0.0137
$ws.Range("N5").ClearContents()

# Sheet 3: AMC Conveyance Solutions Divis
$ws = $wb.Worksheets.Item(3)
$ws.Range("D2").Value = 0.0675
$ws.Range("D3").Value = 0.0675
$ws.Range("D4").Value = 0.0675
$ws.Range("I4").Value = 0.0316
$ws.Range("K4").Value = 0.0072
$ws.Range("M4").Value = 0.0305
$ws.Range("N4").Value = 0.0054
$ws.Range("O4").Value = 0.00964166666666667
$ws.Range("P4").Value = 0.00964166666666667
$ws.Range("Q4").Value = 0.028925
$ws.Range("R4").Value = 0.00964166666666667
$ws.Range("S4").Value = 0.00964166666666667
$ws.Range("T4").Value = 0.00964166666666667
$ws.Range("U4").Value = 0.028925
$ws.Range("V4").Value = 0.1157
$ws.Range("D5").Value = 0.59375
$ws.Range("D6").Value = 0.59375
$ws.Range("D7").Value = 0.59375
$ws.Range("N7").Value = 0.5
$ws.Range("O7").Value = 0.59375
$ws.Range("P7").Value = 0.59375
$ws.Range("Q7").Value = 0.59375
$ws.Range("R7").Value = 0.59375
$ws.Range("S7").Value = 0.59375
$ws.Range("T7").Value = 0.59375
$ws.Range("U7").Value = 0.59375
$ws.Range("V7").Value = 0.59375

# Sheet 5: AMC Linear Motion Division
$ws = $wb.Worksheets.Item(5)
$ws.Range("D2").Value = 0.0469
$ws.Range("D3").Value = 0.0469
$ws.Range("D4").Value = 0.0469
$ws.Range("F4").Value = 0.01
$ws.Range("I4").Value = 0.0222
$ws.Range("M4").Value = 0.0143
$ws.Range("N4").Value = 0.0102
$ws.Range("O4").Value = 0.0067
$ws.Range("P4").Value = 0.0067
$ws.Range("Q4").Value = 0.0201
$ws.Range("R4").Value = 0.0067
$ws.Range("S4").Value = 0.0067
$ws.Range("T4").Value = 0.0067
$ws.Range("U4").Value = 0.0201
$ws.Range("V4").Value = 0.0804
$ws.Range("D5").Value = 0.666666666666667
$ws.Range("D6").Value = 0.666666666666667
$ws.Range("D7").Value = 0.666666666666667
$ws.Range("N7").Value = 0
$ws.Range("O7").Value = 0.666666666666667
$ws.Range("P7").Value = 0.666666666666667
$ws.Range("Q7").Value = 0.666666666666667
$ws.Range("R7").Value = 0.666666666666667
$ws.Range("S7").Value = 0.666666666666667
$ws.Range("T7").Value = 0.666666666666667
$ws.Range("U7").Value = 0.666666666666667
$ws.Range("V7").Value = 0.666666666666667

# Sheet 6: AMC Micro-Motion Division
$ws = $wb.Worksheets.Item(6)
$ws.Range("D2").Value = 0.0369
$ws.Range("D3").Value = 0.0369
$ws.Range("D4").Value = 0.0369
$ws.Range("I4").Value = 0.014
$ws.Range("M4").Value = 0.0172
$ws.Range("N4").Value = 0.0058
$ws.Range("O4").Value = 0.005275
$ws.Range("P4").Value = 0.005275
$ws.Range("Q4").Value = 0.015825
$ws.Range("R4").Value = 0.005275
$ws.Range("S4").Value = 0.005275
$ws.Range("T4").Value = 0.005275
$ws.Range("U4").Value = 0.015825
$ws.Range("V4").Value = 0.0633
$ws.Range("N7").ClearContents()

# Sheet 7: AMC Motion Control Systems Div
$ws = $wb.Worksheets.Item(7)
$ws.Range("D2").Value = 0.0734
$ws.Range("D3").Value = 0.0734
$ws.Range("D4").Value = 0.0734
$ws.Range("F4").Value = 0.0137
$ws.Range("G4").Value = 0.0078
$ws.Range("I4").Value = 0.0279
$ws.Range("J4").Value = 0.0093
$ws.Range("K4").Value = 0.0157
$ws.Range("L4").Value = 0.0127
$ws.Range("M4").Value = 0.0377
$ws.Range("N4").Value = 0.0079
$ws.Range("O4").Value = 0.0104833333333333
$ws.Range("P4").Value = 0.0104833333333333
$ws.Range("Q4").Value = 0.03145
$ws.Range("R4").Value = 0.0104833333333333
$ws.Range("S4").Value = 0.0104833333333333
$ws.Range("T4").Value = 0.0104833333333333
$ws.Range("U4").Value = 0.03145
$ws.Range("V4").Value = 0.1258
$ws.Range("D5").Value = 0.875
$ws.Range("D6").Value = 0.875
$ws.Range("D7").Value = 0.875
$ws.Range("L7").Value = 0.8571
$ws.Range("M7").Value = 0.8667
$ws.Range("N7").Value = 1
$ws.Range("O7").Value = 0.875
$ws.Range("P7").Value = 0.875
$ws.Range("Q7").Value = 0.875
$ws.Range("R7").Value = 0.875
$ws.Range("S7").Value = 0.875
$ws.Range("T7").Value = 0.875
$ws.Range("U7").Value = 0.875
$ws.Range("V7").Value = 0.875

# Sheet 8: AMC Power Management Division
$ws = $wb.Worksheets.Item(8)
$ws.Range("D2").Value = 0.0701
$ws.Range("D3").Value = 0.0701
$ws.Range("D4").Value = 0.0701
$ws.Range("I4").Value = 0.0174
$ws.Range("J4").Value = 0.0115
$ws.Range("K4").Value = 0.0176
$ws.Range("M4").Value = 0.035
$ws.Range("N4").Value = 0.018
$ws.Range("O4").Value = 0.0100166666666667
$ws.Range("P4").Value = 0.0100166666666667
$ws.Range("Q4").Value = 0.03005
$ws.Range("R4").Value = 0.0100166666666667
$ws.Range("S4").Value = 0.0100166666666667
$ws.Range("T4").Value = 0.0100166666666667
$ws.Range("U4").Value = 0.03005
$ws.Range("V4").Value = 0.1202
$ws.Range("N7").ClearContents()

# Sheet 9: AMC Segment Functions
$ws = $wb.Worksheets.Item(9)
$ws.Range("D2").Value = 0.1205
$ws.Range("D3").Value = 0.1205
$ws.Range("D4").Value = 0.1205
$ws.Range("N4").Value = 0
$ws.Range("O4").Value = 0.0172166666666667
$ws.Range("P4").Value = 0.0172166666666667
$ws.Range("Q4").Value = 0.05165
$ws.Range("R4").Value = 0.0172166666666667
$ws.Range("S4").Value = 0.0172166666666667
$ws.Range("T4").Value = 0.0172166666666667
$ws.Range("U4").Value = 0.05165
$ws.Range("V4").Value = 0.2066
$ws.Range("N7").ClearContents()

# Sheet 10: AMC Thomson Linear Motion - Ge
$ws = $wb.Worksheets.Item(10)
$ws.Range("N4").ClearContents()
$ws.Range("O4").Value = 0.071425
$ws.Range("P4").Value = 0.071425
$ws.Range("Q4").Value = 0.214275
$ws.Range("R4").Value = 0.071425
$ws.Range("S4").Value = 0.071425
$ws.Range("T4").Value = 0.071425
$ws.Range("U4").Value = 0.214275
$ws.Range("V4").Value = 0.8571

# Sheet 11: L1_AMC
$ws = $wb.Worksheets.Item(11)
$ws.Range("D2").Value = 0.0596
$ws.Range("D3").Value = 0.0596
$ws.Range("D4").Value = 0.0596
$ws.Range("I4").Value = 0.0264
$ws.Range("K4").Value = 0.0084
$ws.Range("L4").Value = 0.0076
$ws.Range("M4").Value = 0.026
$ws.Range("N4").Value = 0.0072
$ws.Range("O4").Value = 0.00851666666666667
$ws.Range("P4").Value = 0.00851666666666667
$ws.Range("Q4").Value = 0.02555
$ws.Range("R4").Value = 0.00851666666666667
$ws.Range("S4").Value = 0.00851666666666667
$ws.Range("T4").Value = 0.00851666666666667
$ws.Range("U4").Value = 0.02555
$ws.Range("V4").Value = 0.1022
$ws.Range("D5").Value = 0.617021276595745
$ws.Range("D6").Value = 0.617021276595745
$ws.Range("D7").Value = 0.617021276595745
$ws.Range("L7").Value = 0.625
$ws.Range("M7").Value = 0.6552
$ws.Range("N7").Value = 0.3333
$ws.Range("O7").Value = 0.617021276595745
$ws.Range("P7").Value = 0.617021276595745
$ws.Range("Q7").Value = 0.617021276595745
$ws.Range("R7").Value = 0.617021276595745
$ws.Range("S7").Value = 0.617021276595745
$ws.Range("T7").Value = 0.617021276595745
$ws.Range("U7").Value = 0.617021276595745
$ws.Range("V7").Value = 0.617021276595745
